$d = $word.ActiveDocument

# Change 1: "re plus lisable" -> "re, plus lisable" (add a comma)
$d.Content.Find.Execute("re plus lisable", $false, $false, $false, $false, $false,
                         $true, 1, $false, "re, plus lisable", 2)

# Change 2: remove the stray single-space run that trails the "<lb/>" marker
# in "...l'enchasse avecq un<lb/> " (the run containing just " " is deleted
# outright, leaving the "<lb/>" run untouched).
$rng = $d.Content
$rng.Find.Execute("avecq un<lb/> ")
if ($rng.Find.Found) {
    $spaceRange = $d.Range($rng.End - 1, $rng.End)
    $spaceRange.Delete()
}
